$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label: "Aug 2022" -> "Sep 2022" ---
# Plain Value assignment on a General-formatted cell would get
# auto-recognised as a date ("Sep 2022" parses to a date serial), which
# the source file never intended (the cell stores literal text). Force
# text entry by switching the cell to a text number format first, then
# restore the cell's original (General/default) formatting afterwards by
# pasting the format back in from an untouched neighbour cell (C3), so
# the cell's style index is left exactly as it started.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "Sep 2022"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null

# --- Column D figure updates ---
$ws.Range("D6").Value = 7198061.0
$ws.Range("D7").Value = 670013184.0
$ws.Range("D8").Value = 4692951888.0
$ws.Range("D9").Value = -1520457188.0
$ws.Range("D10").Value = 3849705945.0

$ws.Range("D16").Value = 4306137003.0
$ws.Range("D17").Value = -4306137003.0
$ws.Range("D19").Value = 3849705945.0

$ws.Range("D23").Value = 4527770.0

$ws.Range("D27").Value = 137375358.0

$ws.Range("D30").Value = 3712330587.44
$ws.Range("D31").Value = 3712330587.44
$ws.Range("D32").Value = 3849705945.44
